# ---------------------------------------------------------------------------
# Fixed update to excel issue
#   - Rename "Requested quantity" header to "Weekly_PO_Qty" on "Weekly Quantity"
#   - Rename "Requested quantity" header to "Monthly_PO_Qty" on "Monthly Trend"
#   - Add a new "PO Forecast" sheet with forecasted PO quantities (ds,
#     PO_Forecast, yhat_lower, yhat_upper) for 29 weeks
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" headers -----------------------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the existing sheets ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row -----------------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the header look & feel used on the other sheets (bold, bordered,
# centered) by copying the formatting from an existing header cell instead of
# re-deriving it (keeps the workbook's style table untouched).
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$poForecastData = @(
    @(45081.99999999999, 8, -34.22767860221123, 47.08849365216549),
    @(45095.99999999999, 15, -23.89304996555421, 55.28537488960884),
    @(45102.99999999999, 19, -22.30235624637595, 58.56205235191018),
    @(45109.99999999999, 22, -16.83183370124235, 59.93492954460888),
    @(45116.99999999999, 26, -14.86227989690288, 64.28700218924138),
    @(45123.99999999999, 29, -7.712269507211516, 69.88398683057552),
    @(45130.99999999999, 33, -7.26344119988509, 72.27393303701398),
    @(45137.99999999999, 37, -0.743419319508976, 77.0974193462344),
    @(45144.99999999999, 40, 0.4720540518973372, 79.9793016820035),
    @(45151.99999999999, 44, 6.899355544601818, 85.17925151207253),
    @(45158.99999999999, 48, 7.893896243768209, 85.07273180160928),
    @(45165.99999999999, 51, 9.481858131917637, 87.23139972436978),
    @(45172.99999999999, 55, 15.62846885897153, 93.33278973017164),
    @(45179.99999999999, 58, 20.41178236395189, 95.84064683980726),
    @(45186.99999999999, 62, 23.15985434892488, 99.44486035976324),
    @(45193.99999999999, 66, 24.89626785971113, 101.9596047614701),
    @(45200.99999999999, 69, 30.55316037437408, 105.025973578972),
    @(45207.99999999999, 73, 30.96778603359694, 109.4300852833541),
    @(45214.99999999999, 77, 38.34288671723811, 119.016591444511),
    @(45221.99999999999, 80, 42.20055688078967, 118.6771951963084),
    @(45228.99999999999, 84, 45.65787094247745, 121.094613656624),
    @(45235.99999999999, 87, 46.06805521981763, 125.5953680559394),
    @(45242.99999999999, 91, 52.54185587669614, 132.9462316431638),
    @(45249.99999999999, 95, 53.86475167439077, 132.3180414015164),
    @(45256.99999999999, 98, 63.51970583141605, 137.7131000670386),
    @(45263.99999999999, 102, 61.55021313375509, 141.949762721796),
    @(45270.99999999999, 105, 67.38356560663279, 145.7563167954695),
    @(45277.99999999999, 109, 68.31545999753018, 145.6095573764974),
    @(45284.99999999999, 113, 74.82509948334922, 150.7265624985274),
)

# --- 4. Data rows (weekly forecast, row 2 .. row 30) ------------------------
$row = 2
foreach ($r in $poForecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Column A holds dates (serial numbers) — copy the existing date formatting
# used by the other sheets onto the new column instead of re-creating it, so
# no extra number-format / style entries are introduced.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A30").PasteSpecial(-4122)

